$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text formatting so
# that numeric-looking strings (e.g. "1.005", "0.00001095") are not
# reinterpreted by Excel as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.016.81"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.860.24"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "312.43"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "0.5090"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").Value = "0.3840"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "0.08220"
$ws.Range("E9").Value = "  -9.05%  "
$ws.Range("D10").Value = "1.111"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").Value = "41.50"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").Value = "6.211"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").Value = "20.56"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").Value = "1.858.12"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").Value = "7.253"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "0.00001095"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "90.76"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "0.06642"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "17.70"
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "6.020"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").Value = "28.046.73"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "11.08"
$ws.Range("E24").Value = "  -3.37%  "
$ws.Range("D25").Value = "2.242"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.072.16"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.516"
$ws.Range("E27").Value = "  -1.32%  "
$ws.Range("D28").Value = "157.11"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "20.49"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").Value = "124.78"
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("D31").Value = "0.1061"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "1.035"
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("D33").Value = "5.949"
$ws.Range("E33").Value = "  +5.83%  "
$ws.Range("D34").Value = "3.594"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "9.364"
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("D36").Value = "0.06527"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").Value = "0.02417"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("D38").Value = "0.2172"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").Value = "0.6561"
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("D40").Value = "1.198"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").Value = "5.033"
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("D42").Value = "1.221"
$ws.Range("E42").Value = "  -5.63%  "
$ws.Range("D43").Value = "11.18"
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("D44").Value = "0.6148"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("D45").Value = "13.05"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").Value = "1.281"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").Value = "3.652"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "2.013"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "1.206"
$ws.Range("E49").Value = "  -2.55%  "
$ws.Range("D50").Value = "120.10"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").Value = "78.41"
$ws.Range("E51").Value = "  -1.67%  "
